$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

$ws.Range("B2").Value = 285.96523861719356
$ws.Range("C2").Value = 260.03792218044629
$ws.Range("D2").Value = 286.67860187027861
$ws.Range("E2").Value = 255.93217930074152

$ws.Range("B3").Value = 293.78520853744556
$ws.Range("C3").ClearContents() | Out-Null
$ws.Range("D3").Value = 304.39108615447481
$ws.Range("E3").Value = 249.04615181701169

$ws.Range("B1:E3").Select() | Out-Null
